# Update "想去人数" (F column) counts on the 展览, 演出 and 全部类型 sheets.
# These are simple value bumps reflecting freshly scraped attendance counts.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 1103
$ws.Range("F11").Value = 3053
$ws.Range("F12").Value = 559
$ws.Range("F13").Value = 1721
$ws.Range("F17").Value = 1446
$ws.Range("F23").Value = 58
$ws.Range("F24").Value = 4650
$ws.Range("F28").Value = 39
$ws.Range("F29").Value = 85

# 演出 (Performance) sheet
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 20
$ws.Range("F6").Value = 60
$ws.Range("F7").Value = 21
$ws.Range("F9").Value = 48

# 全部类型 (All Types) sheet
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value  = 20
$ws.Range("F9").Value  = 60
$ws.Range("F11").Value = 21
$ws.Range("F13").Value = 48
$ws.Range("F15").Value = 1103
$ws.Range("F22").Value = 3053
$ws.Range("F23").Value = 559
$ws.Range("F24").Value = 1721
$ws.Range("F28").Value = 1446
$ws.Range("F36").Value = 58
$ws.Range("F37").Value = 4650
$ws.Range("F43").Value = 39
$ws.Range("F44").Value = 85
